$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("JdT-TPI_LRD")

# New work-diary entry for the morning: 13/05/2022, "Réalisation", 0.75h,
# bug-fix work ahead of the acceptance tests.
$ws.Range("A49").Copy()
$ws.Range("A50").PasteSpecial(-4122)

$ws.Range("A50").Value = 44694
$ws.Range("B50").Value = "Réalisation"
$ws.Range("C50").Value = 0.75
$ws.Range("D50").Value = "Correction des bugs présents afin de passer les tests d'acceptations"

# The table auto-expands to include the new row.
$tbl = $ws.ListObjects.Item("Tableau1")
$tbl.Resize($ws.Range("A1:F50"))

$ws.Range("F50").Select()
